$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A26:A34").ClearContents()

$ws.Range("D19").Select()
